$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $text) {
    $rng = $ws.Range($cell)
    $origStyle = $rng.Style
    $rng.Value = "'" + $text
    $rng.Style = $origStyle
}

Set-TextValue "D2" "36.542.75"
Set-TextValue "E2" "  -1.85%  "
Set-TextValue "D3" "2.019.24"
Set-TextValue "E3" "  +0.72%  "
Set-TextValue "E4" "  +0.04%  "
Set-TextValue "D5" "234.15"
Set-TextValue "E5" "  -9.57%  "
Set-TextValue "D6" "0.600"
Set-TextValue "E6" "  -2.54%  "
Set-TextValue "E7" "  +0.09%  "
Set-TextValue "D8" "54.88"
Set-TextValue "E8" "  -2.84%  "
Set-TextValue "E9" "  -2.72%  "
Set-TextValue "D10" "57.36"
Set-TextValue "E10" "  +2.33%  "
Set-TextValue "E11" "  -2.93%  "
Set-TextValue "E12" "  -0.72%  "
Set-TextValue "D13" "2.317.21"
Set-TextValue "E13" "  +0.71%  "
Set-TextValue "D14" "14.19"
Set-TextValue "E14" "  -0.09%  "
Set-TextValue "D15" "20.09"
Set-TextValue "E15" "  -6.93%  "
Set-TextValue "D16" "0.763"
Set-TextValue "E16" "  -3.66%  "
Set-TextValue "D17" "5.09"
Set-TextValue "E17" "  -2.17%  "
Set-TextValue "D18" "2.023.25"
Set-TextValue "E18" "  +0.42%  "
Set-TextValue "D19" "36.443.44"
Set-TextValue "E19" "  -2.08%  "
Set-TextValue "D20" "67.67"
Set-TextValue "E20" "  -3.41%  "
Set-TextValue "E21" "  -4.38%  "
Set-TextValue "D22" "5.37"
Set-TextValue "E22" "  +5.53%  "
Set-TextValue "D23" "219.88"
Set-TextValue "E23" "  -6.03%  "
Set-TextValue "E24" "  +0.07%  "
Set-TextValue "E25" "  +1.28%  "
Set-TextValue "E26" "  -6.79%  "
Set-TextValue "D27" "162.91"
Set-TextValue "E27" "  -1.18%  "
Set-TextValue "E28" "  -4.10%  "
Set-TextValue "D29" "0.128"
Set-TextValue "E29" "  -0.04%  "
Set-TextValue "E30" "  +4.19%  "
Set-TextValue "D31" "18.89"
Set-TextValue "E31" "  -3.26%  "
Set-TextValue "E32" "  -1.95%  "
Set-TextValue "D33" "4.35"
Set-TextValue "E33" "  -4.96%  "
Set-TextValue "D34" "0.0603"
Set-TextValue "E34" "  -5.40%  "
Set-TextValue "E35" "  +4.55%  "
Set-TextValue "E36" "  -4.43%  "
Set-TextValue "E37" "  -0.01%  "
Set-TextValue "E38" "  -2.31%  "
Set-TextValue "E39" "  -3.08%  "
Set-TextValue "D40" "5.67"
Set-TextValue "E40" "  +4.30%  "
Set-TextValue "D41" "2.95"
Set-TextValue "E41" "  -2.91%  "
Set-TextValue "D42" "0.0941"
Set-TextValue "E42" "  +1.92%  "
Set-TextValue "D43" "1.457.36"
Set-TextValue "E43" "  +1.36%  "
Set-TextValue "D44" "4.26"
Set-TextValue "E44" "  +42.35%  "
Set-TextValue "E45" "  -3.31%  "
Set-TextValue "D46" "89.90"
Set-TextValue "E46" "  +0.94%  "
Set-TextValue "E47" "  -6.51%  "
Set-TextValue "D48" "15.35"
Set-TextValue "E48" "  -1.49%  "
Set-TextValue "D49" "1.00"
Set-TextValue "E49" "  -1.44%  "
Set-TextValue "E50" "  -1.59%  "
Set-TextValue "D51" "6.85"
Set-TextValue "E51" "  -1.66%  "
